# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 2220
$ws1.Range("F5").Value = 13316
$ws1.Range("F11").Value = 998
$ws1.Range("F12").Value = 13814
$ws1.Range("F13").Value = 14462
$ws1.Range("F25").Value = 5512
$ws1.Range("F26").Value = 942
$ws1.Range("F27").Value = 924
$ws1.Range("F28").Value = 351
$ws1.Range("F30").Value = 108

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 2220
$ws4.Range("F5").Value = 13316
$ws4.Range("F12").Value = 998
$ws4.Range("F13").Value = 13814
$ws4.Range("F14").Value = 14462
$ws4.Range("F26").Value = 5512
$ws4.Range("F27").Value = 942
$ws4.Range("F28").Value = 924
$ws4.Range("F29").Value = 351
$ws4.Range("F31").Value = 108
